# Refresh the cryptocurrency price/volume snapshot to the latest scrape values,
# including fixing the Maker/VeChain row swap, per the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (exactly as authored in the source sheet, incl.
# the deliberate leading/trailing padding spaces on the percentage cells).
$updates = [ordered]@{
    "D2" = "26.147.06"
    "E2" = "  +3.51%  "
    "D3" = "1.604.16"
    "E3" = "  +3.48%  "
    "D5" = "212.93"
    "E5" = "  +3.06%  "
    "E6" = "  -0.29%  "
    "D7" = "0.486"
    "E7" = "  +2.11%  "
    "E8" = "  +2.71%  "
    "E9" = "  +1.31%  "
    "E10" = "  +2.48%  "
    "E11" = "  +4.76%  "
    "D12" = "1.827.71"
    "E12" = "  +3.52%  "
    "D13" = "1.606.15"
    "E13" = "  +3.38%  "
    "E14" = "  +0.82%  "
    "D15" = "0.512"
    "E15" = "  +1.76%  "
    "D16" = "26.147.69"
    "E16" = "  +3.64%  "
    "E17" = "  +3.57%  "
    "D18" = "0.0₃0723"
    "E18" = "  +2.53%  "
    "E19" = "  -0.29%  "
    "D20" = "203.24"
    "E20" = "  +9.86%  "
    "E21" = "  +3.31%  "
    "E22" = "  +0.71%  "
    "D23" = "6.00"
    "E23" = "  +2.89%  "
    "D24" = "1.87"
    "E24" = "  +15.30%  "
    "D25" = "141.32"
    "E25" = "  +1.56%  "
    "E26" = "  -0.27%  "
    "E27" = "  -4.28%  "
    "D28" = "15.17"
    "E28" = "  +2.67%  "
    "E29" = "  +0.93%  "
    "E30" = "  +2.14%  "
    "E31" = "  +2.15%  "
    "E32" = "  +3.22%  "
    "D33" = "2.97"
    "E33" = "  +0.50%  "
    "E34" = "  +1.93%  "
    "E35" = "  +1.05%  "
    "B36" = "Maker"
    "C36" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D36" = "1.124.43"
    "E36" = "  +4.00%  "
    "B37" = "VeChain"
    "C37" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D37" = "0.0164"
    "E37" = "  +10.83%  "
    "E38" = "  -0.30%  "
    "D39" = "0.785"
    "E39" = "  +3.21%  "
    "E40" = "  +2.70%  "
    "E41" = "  +0.13%  "
    "D42" = "0.784"
    "E42" = "  -1.77%  "
    "E43" = "  +2.63%  "
    "D44" = "1.739.77"
    "E44" = "  +3.55%  "
    "D45" = "92.70"
    "E45" = "  +0.22%  "
    "E46" = "  +4.87%  "
    "D47" = "53.46"
    "E47" = "  +2.55%  "
    "D48" = "0.0504"
    "E48" = "  +0.49%  "
    "E49" = "  +1.16%  "
    "E50" = "  -0.10%  "
    "D51" = "0.0₇0927"
    "E51" = "  -16.60%  "
}

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $range = $ws.Range($cell)

    # Plain-decimal-looking values in the Price column ("6.00", "0.0164", ...) need
    # to be forced to Text first, otherwise Excel re-types them as numbers on entry
    # and trailing zeros / the exact textual form would be lost (e.g. "6.00" -> 6).
    if ($cell.StartsWith("D") -and ($value.Trim() -match "^-?[0-9]+(\.[0-9]+)?$")) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
